$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "M1"
$ws.Cells.Item(2, 2).Value = "Ccl24"
$ws.Cells.Item(2, 3).Value = "Ccr2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.852535666666667
$ws.Cells.Item(2, 8).Value = 5.557607
$ws.Cells.Item(2, 9).Value = 0.315753910656914
$ws.Cells.Item(2, 10).Value = 0.315753910656914
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.08299899999999999
$ws.Cells.Item(2, 14).Value = 0.165998
$ws.Cells.Item(2, 15).Value = 0.0002641221862538864
$ws.Cells.Item(2, 16).Value = 0.0001761022045687984
$ws.Cells.Item(2, 17).Value = 0.1537586077976666
$ws.Cells.Item(2, 18).Value = 0.9225516467859999
$ws.Cells.Item(2, 19).Value = 0.00008339761320091846
$ws.Cells.Item(2, 20).Value = 0.00005560495976790196

# Row 3
$ws.Cells.Item(3, 1).Value = "M1"
$ws.Cells.Item(3, 2).Value = "Ccl24"
$ws.Cells.Item(3, 3).Value = "Ccr2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.852535666666667
$ws.Cells.Item(3, 8).Value = 5.557607
$ws.Cells.Item(3, 9).Value = 0.315753910656914
$ws.Cells.Item(3, 10).Value = 0.315753910656914
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.03385166666666667
$ws.Cells.Item(3, 14).Value = 0.101555
$ws.Cells.Item(3, 15).Value = 0.0001077239027980797
$ws.Cells.Item(3, 16).Value = 0.000107736595531177
$ws.Cells.Item(3, 17).Value = 0.06271141987611112
$ws.Cells.Item(3, 18).Value = 0.564402778885
$ws.Cells.Item(3, 19).Value = 0.00003401424357971894
$ws.Cells.Item(3, 20).Value = 0.00003401825135983135

# Row 4
$ws.Cells.Item(4, 1).Value = "M1"
$ws.Cells.Item(4, 2).Value = "Ccl24"
$ws.Cells.Item(4, 3).Value = "Ccr2"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.852535666666667
$ws.Cells.Item(4, 8).Value = 5.557607
$ws.Cells.Item(4, 9).Value = 0.315753910656914
$ws.Cells.Item(4, 10).Value = 0.315753910656914
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 136.0989836666667
$ws.Cells.Item(4, 14).Value = 408.296951
$ws.Cells.Item(4, 15).Value = 0.4330987254421377
$ws.Cells.Item(4, 16).Value = 0.4331497559598229
$ws.Cells.Item(4, 17).Value = 252.1282214395841
$ws.Cells.Item(4, 18).Value = 2269.153992956257
$ws.Cells.Item(4, 19).Value = 0.1367526162588801
$ws.Cells.Item(4, 20).Value = 0.136768729344402

# Row 5
$ws.Cells.Item(5, 1).Value = "M1"
$ws.Cells.Item(5, 2).Value = "Ccl24"
$ws.Cells.Item(5, 3).Value = "Ccr2"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.852535666666667
$ws.Cells.Item(5, 8).Value = 5.557607
$ws.Cells.Item(5, 9).Value = 0.315753910656914
$ws.Cells.Item(5, 10).Value = 0.315753910656914
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 177.6775513333333
$ws.Cells.Item(5, 14).Value = 533.032654
$ws.Cells.Item(5, 15).Value = 0.5654114303352707
$ws.Cells.Item(5, 16).Value = 0.5654780507991515
$ws.Cells.Item(5, 17).Value = 329.1540010109975
$ws.Cells.Item(5, 18).Value = 2962.386009098978
$ws.Cells.Item(5, 19).Value = 0.178530870258481
$ws.Cells.Item(5, 20).Value = 0.1785519059304811

# Row 6
$ws.Cells.Item(6, 1).Value = "M1"
$ws.Cells.Item(6, 2).Value = "Ccl24"
$ws.Cells.Item(6, 3).Value = "Ccr2"
$ws.Cells.Item(6, 4).Value = "Neutro"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.852535666666667
$ws.Cells.Item(6, 8).Value = 5.557607
$ws.Cells.Item(6, 9).Value = 0.315753910656914
$ws.Cells.Item(6, 10).Value = 0.315753910656914
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.323258
$ws.Cells.Item(6, 14).Value = 0.9697739999999999
$ws.Cells.Item(6, 15).Value = 0.001028682389957214
$ws.Cells.Item(6, 16).Value = 0.001028803596028277
$ws.Cells.Item(6, 17).Value = 0.5988469745353333
$ws.Cells.Item(6, 18).Value = 5.389622770818
$ws.Cells.Item(6, 19).Value = 0.000324810487452891
$ws.Cells.Item(6, 20).Value = 0.0003248487587438244

# Row 7
$ws.Cells.Item(7, 1).Value = "M1"
$ws.Cells.Item(7, 2).Value = "Ccl24"
$ws.Cells.Item(7, 3).Value = "Ccr2"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.852535666666667
$ws.Cells.Item(7, 8).Value = 5.557607
$ws.Cells.Item(7, 9).Value = 0.315753910656914
$ws.Cells.Item(7, 10).Value = 0.315753910656914
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.5
$ws.Cells.Item(7, 13).Value = 0.028067
$ws.Cells.Item(7, 14).Value = 0.056134
$ws.Cells.Item(7, 15).Value = 0.0000893157435823062
$ws.Cells.Item(7, 16).Value = 0.00005955084489731763
$ws.Cells.Item(7, 17).Value = 0.05199511855633333
$ws.Cells.Item(7, 18).Value = 0.311970711338
$ws.Cells.Item(7, 19).Value = 0.00002820179531934335
$ws.Cells.Item(7, 20).Value = 0.00001880341215925137

# Row 8
$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Ccl24"
$ws.Cells.Item(8, 3).Value = "Ccr2"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 4.014488
$ws.Cells.Item(8, 8).Value = 12.043464
$ws.Cells.Item(8, 9).Value = 0.6842460893430861
$ws.Cells.Item(8, 10).Value = 0.684246089343086
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.08299899999999999
$ws.Cells.Item(8, 14).Value = 0.165998
$ws.Cells.Item(8, 15).Value = 0.0002641221862538864
$ws.Cells.Item(8, 16).Value = 0.0001761022045687984
$ws.Cells.Item(8, 17).Value = 0.333198489512
$ws.Cells.Item(8, 18).Value = 1.999190937072
$ws.Cells.Item(8, 19).Value = 0.000180724573052968
$ws.Cells.Item(8, 20).Value = 0.0001204972448008964

# Row 9
$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Ccl24"
$ws.Cells.Item(9, 3).Value = "Ccr2"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 4.014488
$ws.Cells.Item(9, 8).Value = 12.043464
$ws.Cells.Item(9, 9).Value = 0.6842460893430861
$ws.Cells.Item(9, 10).Value = 0.684246089343086
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.03385166666666667
$ws.Cells.Item(9, 14).Value = 0.101555
$ws.Cells.Item(9, 15).Value = 0.0001077239027980797
$ws.Cells.Item(9, 16).Value = 0.000107736595531177
$ws.Cells.Item(9, 17).Value = 0.1358971096133333
$ws.Cells.Item(9, 18).Value = 1.22307398652
$ws.Cells.Item(9, 19).Value = 0.00007370965921836075
$ws.Cells.Item(9, 20).Value = 0.00007371834417134568

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Ccl24"
$ws.Cells.Item(10, 3).Value = "Ccr2"
$ws.Cells.Item(10, 4).Value = "M1"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4.014488
$ws.Cells.Item(10, 8).Value = 12.043464
$ws.Cells.Item(10, 9).Value = 0.6842460893430861
$ws.Cells.Item(10, 10).Value = 0.684246089343086
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 136.0989836666667
$ws.Cells.Item(10, 14).Value = 408.296951
$ws.Cells.Item(10, 15).Value = 0.4330987254421377
$ws.Cells.Item(10, 16).Value = 0.4331497559598229
$ws.Cells.Item(10, 17).Value = 546.3677367420294
$ws.Cells.Item(10, 18).Value = 4917.309630678265
$ws.Cells.Item(10, 19).Value = 0.2963461091832577
$ws.Cells.Item(10, 20).Value = 0.2963810266154208

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Ccl24"
$ws.Cells.Item(11, 3).Value = "Ccr2"
$ws.Cells.Item(11, 4).Value = "M2"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 4.014488
$ws.Cells.Item(11, 8).Value = 12.043464
$ws.Cells.Item(11, 9).Value = 0.6842460893430861
$ws.Cells.Item(11, 10).Value = 0.684246089343086
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 177.6775513333333
$ws.Cells.Item(11, 14).Value = 533.032654
$ws.Cells.Item(11, 15).Value = 0.5654114303352707
$ws.Cells.Item(11, 16).Value = 0.5654780507991515
$ws.Cells.Item(11, 17).Value = 713.2843976970506
$ws.Cells.Item(11, 18).Value = 6419.559579273456
$ws.Cells.Item(11, 19).Value = 0.3868805600767897
$ws.Cells.Item(11, 20).Value = 0.3869261448686703

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Ccl24"
$ws.Cells.Item(12, 3).Value = "Ccr2"
$ws.Cells.Item(12, 4).Value = "Neutro"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 4.014488
$ws.Cells.Item(12, 8).Value = 12.043464
$ws.Cells.Item(12, 9).Value = 0.6842460893430861
$ws.Cells.Item(12, 10).Value = 0.684246089343086
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.323258
$ws.Cells.Item(12, 14).Value = 0.9697739999999999
$ws.Cells.Item(12, 15).Value = 0.001028682389957214
$ws.Cells.Item(12, 16).Value = 0.001028803596028277
$ws.Cells.Item(12, 17).Value = 1.297715361904
$ws.Cells.Item(12, 18).Value = 11.679438257136
$ws.Cells.Item(12, 19).Value = 0.0007038719025043233
$ws.Cells.Item(12, 20).Value = 0.0007039548372844524

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Ccl24"
$ws.Cells.Item(13, 3).Value = "Ccr2"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 4.014488
$ws.Cells.Item(13, 8).Value = 12.043464
$ws.Cells.Item(13, 9).Value = 0.6842460893430861
$ws.Cells.Item(13, 10).Value = 0.684246089343086
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.5
$ws.Cells.Item(13, 13).Value = 0.028067
$ws.Cells.Item(13, 14).Value = 0.056134
$ws.Cells.Item(13, 15).Value = 0.0000893157435823062
$ws.Cells.Item(13, 16).Value = 0.00005955084489731763
$ws.Cells.Item(13, 17).Value = 0.112674634696
$ws.Cells.Item(13, 18).Value = 0.6760478081760001
$ws.Cells.Item(13, 19).Value = 0.00006111394826296284
$ws.Cells.Item(13, 20).Value = 0.00004074743273806625
